$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update status of row 18 (Mouse Pointer / curser issue) from "Not Started" to "Complete"
$ws.Range("F18").Value = "Complete"

# Set Started / Completed dates for row 18 (reuse the existing date-formatted
# style from another row so we don't mint a brand-new number format)
$ws.Range("G3").Copy() | Out-Null
$ws.Range("G18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("G18").Value = 42984

$ws.Range("H3").Copy() | Out-Null
$ws.Range("H18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("H18").Value = 42984

$excel.CutCopyMode = 0

# Update the view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E15").Select()
